$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "Play Boots of Luck Slot Game for Free - Review" "Play Boots of Luck Slot Game for Free"
Replace-Text "Colorful and straightforward design" "Colorful and straightforward gameplay"
Replace-Text "Utilizes a common reel system with 20 paylines" "Bright graphics and sound for a fun experience"
Replace-Text "Features a wild and scatter symbol" "Irish-themed design adds a unique touch"
Replace-Text "Similar games available for players to enjoy" "Free spins feature for additional chances to win"
Replace-Text "Only five free spins rewarded during bonus mode" "Limited number of paylines compared to some other slot games"
Replace-Text "Notably similar to other Irish-themed slot games" "Not as many bonus features as some players may prefer"
Replace-Text "Read our review of Boots of Luck, an Irish-themed slot game with a 5x3 reel system and 20 paylines. Play now for free and enjoy the colorful design and wild and scatter symbols." "Read our review of Boots of Luck, an Irish-themed slot game by Betixon. Play for free and enjoy colorful graphics and chances to win!"
